$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.776953935623169
$ws.Range("B1").Value = 3.834307670593262
$ws.Range("C1").Value = 1.297136545181274
$ws.Range("D1").Value = 0.8588575720787048
$ws.Range("E1").Value = 0.4611234068870544
